# Apply cryptocurrency symbol-list refresh as captured in the commit
# "Updated symbol list on Sun Jan 15 05:53:17 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell and its new text value.
$updates = @(
    @{ Cell = "D2"; Value = "298.58" },
    @{ Cell = "E2"; Value = "-2.71%" },
    @{ Cell = "D3"; Value = "31.78" },
    @{ Cell = "E3"; Value = "-1.51%" },
    @{ Cell = "D4"; Value = "5.101" },
    @{ Cell = "E4"; Value = "-4.36%" },
    @{ Cell = "D5"; Value = "0.07535" },
    @{ Cell = "E5"; Value = "1.40%" },
    @{ Cell = "B6"; Value = "KuCoinToken" },
    @{ Cell = "C6"; Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs" },
    @{ Cell = "D6"; Value = "7.744" },
    @{ Cell = "E6"; Value = "-0.44%" },
    @{ Cell = "B7"; Value = "FTXToken" },
    @{ Cell = "C7"; Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt" },
    @{ Cell = "D7"; Value = "1.763" },
    @{ Cell = "E7"; Value = "12.58%" },
    @{ Cell = "D8"; Value = "3.793" },
    @{ Cell = "E8"; Value = "2.76%" },
    @{ Cell = "D9"; Value = "0.9282" },
    @{ Cell = "E9"; Value = "1.71%" },
    @{ Cell = "D10"; Value = "0.1709" },
    @{ Cell = "E10"; Value = "2.37%" },
    @{ Cell = "D11"; Value = "0.07260" },
    @{ Cell = "E11"; Value = "-4.60%" },
    @{ Cell = "D12"; Value = "0.07948" },
    @{ Cell = "E12"; Value = "-1.57%" },
    @{ Cell = "D13"; Value = "0.03059" },
    @{ Cell = "E13"; Value = "1.62%" },
    @{ Cell = "D14"; Value = "0.09886" },
    @{ Cell = "E14"; Value = "0.36%" },
    @{ Cell = "D15"; Value = "0.001492" },
    @{ Cell = "E15"; Value = "-3.54%" },
    @{ Cell = "D16"; Value = "0.006518" },
    @{ Cell = "E16"; Value = "4.53%" },
    @{ Cell = "D17"; Value = "3.450" },
    @{ Cell = "D18"; Value = "2.221" },
    @{ Cell = "E18"; Value = "-0.71%" },
    @{ Cell = "E19"; Value = "0.50%" },
    @{ Cell = "E20"; Value = "-0.76%" },
    @{ Cell = "D21"; Value = "4.560" },
    @{ Cell = "E21"; Value = "8.19%" },
    @{ Cell = "D22"; Value = "0.04651" },
    @{ Cell = "E22"; Value = "2.11%" },
    @{ Cell = "E23"; Value = "-4.85%" },
    @{ Cell = "D24"; Value = "0.001217" },
    @{ Cell = "E24"; Value = "0.14%" },
    @{ Cell = "D25"; Value = "0.004423" },
    @{ Cell = "E25"; Value = "-1.71%" },
    @{ Cell = "D26"; Value = "0.0001399" },
    @{ Cell = "E27"; Value = "6.76%" },
    @{ Cell = "E39"; Value = "-0.81%" },
    @{ Cell = "D40"; Value = "0.04555" },
    @{ Cell = "E40"; Value = "1.03%" },
    @{ Cell = "D41"; Value = "0.007056" },
    @{ Cell = "E41"; Value = "-2.00%" },
    @{ Cell = "E42"; Value = "-2.49%" },
    @{ Cell = "D44"; Value = "0.01281" },
    @{ Cell = "E44"; Value = "-6.60%" },
    @{ Cell = "D45"; Value = "0.00006049" },
    @{ Cell = "E45"; Value = "-1.13%" },
    @{ Cell = "D46"; Value = "1.930" },
    @{ Cell = "E46"; Value = "1.97%" },
    @{ Cell = "E47"; Value = "-0.20%" }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    # Columns D (Price) and E (Volume) hold numeric-looking text (e.g. "298.58",
    # "-2.71%"); force a Text number format first so Excel stores the new value
    # as a string instead of re-parsing it into a number/percentage.
    if ($u.Cell -match "^[DE]\d+$") {
        $range.NumberFormat = "@"
    }
    $range.Value = $u.Value
}
